{"js": "// Append new rows to the second table (the session log table with\n// Date / Time / Duration / Role / Role columns).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[tables.items.length - 1];\n\nconst newRows = [\n  [\"07/2/23\", \"12:30\", \"1h \", \"Observer\", \"Driver\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n];\n\ntable.addRows(\"End\", newRows.length, newRows);\nawait context.sync();\n", "ps1": "# Append new rows to the last table (the session log table with\n# Date / Time / Duration / Role / Role columns).\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item($d.Tables.Count)\n\n$newRow = $tbl.Rows.Add()\n$newRow.Cells.Item(1).Range.Text = \"07/2/23\"\n$newRow.Cells.Item(2).Range.Text = \"12:30\"\n$newRow.Cells.Item(3).Range.Text = \"1h \"\n$newRow.Cells.Item(4).Range.Text = \"Observer\"\n$newRow.Cells.Item(5).Range.Text = \"Driver\"\n\nfor ($i = 0; $i -lt 7; $i++) {\n    $tbl.Rows.Add() | Out-Null\n}\n"}
